$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLAN CAPITAL")

# Independent single-cell corrections in column L (week of 16/01/2026),
# unrelated to the roster change below.
$ws.Range("L11").Value = 6
$ws.Range("L13").Value = 6
$ws.Range("L17").Value = 6
$ws.Range("L18").Value = 2
$ws.Range("L20").Value = 5
$ws.Range("L24").Value = 6
$ws.Range("L30").Value = 6
$ws.Range("L31").Value = 6

# Remove the "Punpie69" roster row (row 33); everything below shifts up
# by one row. This also pulls the closing thick-bottom border row back up
# to where the data ends, so we re-insert a blank spacer row and wipe it
# so the border row keeps its own row number and formatting (it does not
# belong to the roster table).
$ws.Rows.Item(33).Delete()
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Clear()
